# Insert a new row at position 110 (this shifts the former rows 110..185
# down to 111..186, matching the rest of the diff) and populate the new
# row with the new data record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 110, pushing every row
# from 110 downward one position further down (old 110 -> 111, ... old
# 185 -> 186).
$ws.Rows("110:110").Insert()

# Populate the newly inserted row 110 with the new record's values.
$ws.Range("A110").Value = 4
$ws.Range("B110").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C110").Value = "Los Lagos"
$ws.Range("D110").Value = 44729
$ws.Range("E110").Value = 10
$ws.Range("F110").Value = 100112009
$ws.Range("G110").Value = "Acelga"
$ws.Range("H110").Value = "Sin especificar"
$ws.Range("I110").Value = "Primera"
$ws.Range("J110").Value = 90
$ws.Range("K110").Value = 12000
$ws.Range("L110").Value = 12000
$ws.Range("M110").Value = 12000
$ws.Range("N110").Value = "`$/docena de atados (12 kilos)"
$ws.Range("O110").Value = "Región de La Araucanía"
$ws.Range("P110").Value = 1000
$ws.Range("Q110").Value = 12
$ws.Range("R110").Value = "Hortaliza"
